$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Arabic spelling for "Wadi Ibn Hashbal" (row 167, column C):
# "وادى بن مشيل" -> "وادى بن هشبل"
$ws.Cells.Item(167, 3).Value = "وادى بن هشبل"

# Copy the formatting (border style) of the last existing data row down onto
# the two new rows before filling in their values.
$ws.Range("A174:G174").Copy()
$ws.Range("A175:G176").PasteSpecial(-4122)

# Append two new city rows at the end of the table
$ws.Cells.Item(175, 1).Value = "Hubuna"
$ws.Cells.Item(175, 2).Value = "Hubuna"
$ws.Cells.Item(175, 3).Value = "حبونا"
$ws.Cells.Item(175, 4).Value = 17.839884999999999
$ws.Cells.Item(175, 5).Value = 44.023803999999998
$ws.Cells.Item(175, 6).Value = "منطقة نجران"
$ws.Cells.Item(175, 7).Value = "جنوب المملكة"

$ws.Cells.Item(176, 1).Value = "Tabalah"
$ws.Cells.Item(176, 2).Value = "Tabalah"
$ws.Cells.Item(176, 3).Value = "تبالة"
$ws.Cells.Item(176, 4).Value = 19.996976
$ws.Cells.Item(176, 5).Value = 42.226551999999998
$ws.Cells.Item(176, 6).Value = "منطقة عسير"
$ws.Cells.Item(176, 7).Value = "جنوب المملكة"

# Update the worksheet's used-range selection to include the newly added rows
[void]$ws.Range("A1:G176").Select()
